$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.095.77'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.344.57'
$ws.Range("E3").Value = '  -3.28%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.16'
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '85.24'
$ws.Range("E6").Value = '  -4.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.05'
$ws.Range("E11").Value = '  -5.76%  '
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.705.92'
$ws.Range("E13").Value = '  -3.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.40'
$ws.Range("E14").Value = '  -4.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.75'
$ws.Range("E15").Value = '  -6.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.373.06'
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.759'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '40.061.33'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0901'
$ws.Range("E19").Value = '  -2.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.09'
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.08'
$ws.Range("E21").Value = '  -4.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.66'
$ws.Range("E22").Value = '  -3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.97'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.55'
$ws.Range("E24").Value = '  -5.15%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").Value = '  -3.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.76'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("E28").Value = '  -3.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.27'
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.59'
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.87'
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("E34").Value = '  -2.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0716'
$ws.Range("E35").Value = '  -3.96%  '
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.80'
$ws.Range("E37").Value = '  -6.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0984'
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("E39").Value = '  -4.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '15.53'
$ws.Range("E40").Value = '  -6.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.85'
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.964.74'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("E44").Value = '  -3.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.58'
$ws.Range("E45").Value = '  -6.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.41'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("E47").Value = '  -6.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.567.57'
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.39'
$ws.Range("E50").Value = '  -4.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.18'
$ws.Range("E51").Value = '  -4.06%  '
